# A new daily price record was inserted as row 357 (Provincia de Cautín,
# fecha 2022-07-12 / serial 44754). Every existing record from row 357
# through row 454 shifts down by one row, with the former last record
# (row 454) ending up at the new row 455.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 357:454 down to 358:455, leaving a blank row 357 behind
# (inherits the date-formatted style from the old row 357, same as Excel
# does for a native "Insert" on a full row).
$ws.Rows.Item(357).Insert()

# Fill in the new record in the now-empty row 357.
$ws.Cells.Item(357, 1).Value = 10
$ws.Cells.Item(357, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(357, 3).Value = "La Araucanía"
$ws.Cells.Item(357, 4).Value = 44754
$ws.Cells.Item(357, 5).Value = 9
$ws.Cells.Item(357, 6).Value = 100112023
$ws.Cells.Item(357, 7).Value = "Brócoli"
$ws.Cells.Item(357, 8).Value = "Sin especificar"
$ws.Cells.Item(357, 9).Value = "Primera"
$ws.Cells.Item(357, 10).Value = 500
$ws.Cells.Item(357, 11).Value = 1200
$ws.Cells.Item(357, 12).Value = 1200
$ws.Cells.Item(357, 13).Value = 1200
$ws.Cells.Item(357, 14).Value = "`$/unidad"
$ws.Cells.Item(357, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(357, 16).Value = 1200
$ws.Cells.Item(357, 17).Value = 1
$ws.Cells.Item(357, 18).Value = "Hortaliza"
